$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1927272727272727
$ws.Range("C2").Value = 0.5527272727272727
$ws.Range("J2").Value = 0.01454545454545455
$ws.Range("P2").Value = 0.1272727272727273
$ws.Range("S2").Value = 0.1127272727272727
$ws.Range("B3").Value = 0.0131578947368421
$ws.Range("C3").Value = 0.01973684210526316
$ws.Range("J3").Value = 0.01973684210526316
$ws.Range("P3").Value = 0.743421052631579
$ws.Range("S3").Value = 0.2039473684210526
$ws.Range("J4").Value = 0.09375
$ws.Range("P4").Value = 0.65625
$ws.Range("S4").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.06074766355140187
$ws.Range("F6").Value = 0.08411214953271028
$ws.Range("J6").Value = 0.2429906542056075
$ws.Range("O6").Value = 0.0514018691588785
$ws.Range("Q6").Value = 0.1355140186915888
$ws.Range("R6").Value = 0.06074766355140187
$ws.Range("S6").Value = 0.3644859813084112
$ws.Range("B7").Value = 0.1029411764705882
$ws.Range("D7").Value = 0.01470588235294118
$ws.Range("F7").Value = 0.06862745098039216
$ws.Range("J7").Value = 0.1813725490196078
$ws.Range("O7").Value = 0.02450980392156863
$ws.Range("Q7").Value = 0.1519607843137255
$ws.Range("R7").Value = 0.08333333333333333
$ws.Range("S7").Value = 0.3725490196078431
$ws.Range("B8").Value = 0.05720823798627003
$ws.Range("D8").Value = 0.02288329519450801
$ws.Range("E8").Value = 0.004576659038901602
$ws.Range("F8").Value = 0.08466819221967964
$ws.Range("J8").Value = 0.1075514874141876
$ws.Range("O8").Value = 0.02517162471395881
$ws.Range("Q8").Value = 0.1716247139588101
$ws.Range("R8").Value = 0.08466819221967964
$ws.Range("S8").Value = 0.4416475972540046
$ws.Range("B9").Value = 0.1044776119402985
$ws.Range("D9").Value = 0.009950248756218905
$ws.Range("F9").Value = 0.07960199004975124
$ws.Range("J9").Value = 0.1641791044776119
$ws.Range("O9").Value = 0.02985074626865672
$ws.Range("Q9").Value = 0.1343283582089552
$ws.Range("R9").Value = 0.1442786069651741
$ws.Range("S9").Value = 0.3333333333333333
$ws.Range("B10").Value = 0.1089108910891089
$ws.Range("D10").Value = 0.01485148514851485
$ws.Range("F10").Value = 0.06023102310231023
$ws.Range("J10").Value = 0.132013201320132
$ws.Range("O10").Value = 0.0132013201320132
$ws.Range("Q10").Value = 0.1864686468646865
$ws.Range("R10").Value = 0.1056105610561056
$ws.Range("S10").Value = 0.3787128712871287
$ws.Range("G11").Value = 0.124223602484472
$ws.Range("J11").Value = 0.1118012422360248
$ws.Range("K11").Value = 0.1739130434782609
$ws.Range("L11").Value = 0.5838509316770186
$ws.Range("S11").Value = 0.006211180124223602
$ws.Range("G12").Value = 0.7631578947368421
$ws.Range("J12").Value = 0.1894736842105263
$ws.Range("K12").Value = 0.005263157894736842
$ws.Range("L12").Value = 0.01578947368421053
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("F13").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.5476190476190477
$ws.Range("J13").Value = 0.4047619047619048
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("F15").Value = 0.02439024390243903
$ws.Range("H15").Value = 0.1951219512195122
$ws.Range("I15").Value = 0.07804878048780488
$ws.Range("J15").Value = 0.2341463414634146
$ws.Range("K15").Value = 0.09268292682926829
$ws.Range("M15").Value = 0.00975609756097561
$ws.Range("O15").Value = 0.0975609756097561
$ws.Range("S15").Value = 0.2682926829268293
$ws.Range("F16").Value = 0.01212121212121212
$ws.Range("H16").Value = 0.1454545454545454
$ws.Range("I16").Value = 0.09696969696969697
$ws.Range("J16").Value = 0.4242424242424243
$ws.Range("K16").Value = 0.103030303030303
$ws.Range("M16").Value = 0.02424242424242424
$ws.Range("O16").Value = 0.02424242424242424
$ws.Range("S16").Value = 0.1696969696969697
$ws.Range("F17").Value = 0.01799485861182519
$ws.Range("H17").Value = 0.1773778920308483
$ws.Range("I17").Value = 0.09768637532133675
$ws.Range("J17").Value = 0.4010282776349615
$ws.Range("K17").Value = 0.1182519280205656
$ws.Range("M17").Value = 0.007712082262210797
$ws.Range("O17").Value = 0.05912596401028278
$ws.Range("S17").Value = 0.1208226221079692
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.2142857142857143
$ws.Range("I18").Value = 0.1205357142857143
$ws.Range("J18").Value = 0.3705357142857143
$ws.Range("K18").Value = 0.1071428571428571
$ws.Range("M18").Value = 0.008928571428571428
$ws.Range("O18").Value = 0.02678571428571428
$ws.Range("S18").Value = 0.1205357142857143
$ws.Range("F19").Value = 0.01189532117367169
$ws.Range("H19").Value = 0.2038065027755749
$ws.Range("I19").Value = 0.08009516256938938
$ws.Range("J19").Value = 0.3584456780333069
$ws.Range("K19").Value = 0.1245043616177637
$ws.Range("M19").Value = 0.02696272799365583
$ws.Range("N19").Value = 0.002379064234734338
$ws.Range("O19").Value = 0.06106264869151467
$ws.Range("S19").Value = 0.1308485329103886
